# Auto-generated: applies literal value refresh to H:N columns
# across multiple sheets, per the commit's scheduled-runner data update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 45
$ws.Range("I5").Value = 45
$ws.Range("K5").Value = 45
$ws.Range("M5").Value = 70

$ws.Range("H12").Value = 547.7143
$ws.Range("J12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("N12").Value = -1340

$ws.Range("H28").Value = 750
$ws.Range("I28").Value = 634.4761999999999
$ws.Range("J28").Value = 1154.3334
$ws.Range("K28").Value = 634.4761999999999
$ws.Range("L28").Value = 1154.3334
$ws.Range("M28").Value = -149.4761999999999
$ws.Range("N28").Value = -2124.3334

$ws.Range("H42").Value = 96.5
$ws.Range("J42").Value = 300
$ws.Range("L42").Value = 900
$ws.Range("N42").Value = -1360

$ws.Range("H51").Value = 5197.2
$ws.Range("I51").Value = 4995
$ws.Range("K51").Value = 4995
$ws.Range("M51").Value = -4511

$ws.Range("H61").Value = 145
$ws.Range("I61").Value = 145
$ws.Range("K61").Value = 435
$ws.Range("M61").Value = -263

$ws.Range("H104").Value = 209.16667
$ws.Range("I104").Value = 209.16667
$ws.Range("K104").Value = 627.50001
$ws.Range("M104").Value = 1119.49999

$ws.Range("H111").Value = 500
$ws.Range("I111").Value = 400
$ws.Range("K111").Value = 1200
$ws.Range("M111").Value = 1867

$ws.Range("H115").Value = 192.5
$ws.Range("I115").Value = 192.5
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 577.5
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 989.5
$ws.Range("N115").ClearContents()

$ws.Range("H116").Value = 24999.4
$ws.Range("I116").Value = 6665.6665
$ws.Range("K116").Value = 6665.6665
$ws.Range("M116").Value = -3223.6665

$ws.Range("H131").Value = 2999.5
$ws.Range("I131").Value = 2999.5
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 8998.5
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -3958.5
$ws.Range("N131").ClearContents()

$ws.Range("H137").Value = 1859.8125
$ws.Range("I137").Value = 1554.0714
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 4662.2142
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -2112.2142
$ws.Range("N137").Value = -17100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2086.611
$ws.Range("I74").Value = 2086.611
$ws.Range("K74").Value = 2086.611
$ws.Range("M74").Value = -1212.611

$ws.Range("H77").Value = 2086.611
$ws.Range("I77").Value = 2086.611
$ws.Range("K77").Value = 10433.055
$ws.Range("M77").Value = -6065.055

$ws.Range("H110").Value = 3165.4546
$ws.Range("I110").Value = 2053.3333
$ws.Range("K110").Value = 2053.3333
$ws.Range("M110").Value = -8.333299999999781

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 13112.621
$ws.Range("I134").Value = 10438.833
$ws.Range("K134").Value = 31316.499
$ws.Range("M134").Value = -28781.499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 42.18182
$ws.Range("I7").Value = 45.25
$ws.Range("J7").Value = 34
$ws.Range("K7").Value = 45.25
$ws.Range("L7").Value = 34
$ws.Range("M7").Value = 67.75
$ws.Range("N7").Value = -260

$ws.Range("H16").Value = 8458
$ws.Range("I16").Value = 6149.6
$ws.Range("K16").Value = 6149.6
$ws.Range("M16").Value = -5862.6

$ws.Range("H31").Value = 4918.4
$ws.Range("I31").Value = 5992
$ws.Range("J31").Value = 4753.231
$ws.Range("K31").Value = 5992
$ws.Range("L31").Value = 4753.231
$ws.Range("M31").Value = -5697
$ws.Range("N31").Value = -5343.231

$ws.Range("H34").Value = 4918.4
$ws.Range("I34").Value = 5992
$ws.Range("J34").Value = 4753.231
$ws.Range("K34").Value = 5992
$ws.Range("L34").Value = 4753.231
$ws.Range("M34").Value = -5790
$ws.Range("N34").Value = -5157.231

$ws.Range("H62").Value = 9849.9
$ws.Range("I62").Value = 10266.5
$ws.Range("J62").Value = 9225
$ws.Range("K62").Value = 10266.5
$ws.Range("L62").Value = 9225
$ws.Range("M62").Value = -9642.5
$ws.Range("N62").Value = -10473

$ws.Range("H65").Value = 9849.9
$ws.Range("I65").Value = 10266.5
$ws.Range("J65").Value = 9225
$ws.Range("K65").Value = 51332.5
$ws.Range("L65").Value = 46125
$ws.Range("M65").Value = -48212.5
$ws.Range("N65").Value = -52365

$ws.Range("H99").Value = 7892.857
$ws.Range("I99").Value = 8416.666999999999
$ws.Range("K99").Value = 8416.666999999999
$ws.Range("M99").Value = -6918.666999999999

$ws.Range("H113").Value = 8458
$ws.Range("I113").Value = 6149.6
$ws.Range("K113").Value = 6149.6
$ws.Range("M113").Value = -3979.6

$ws.Range("H126").Value = 7892.857
$ws.Range("I126").Value = 8416.666999999999
$ws.Range("K126").Value = 25250.001
$ws.Range("M126").Value = -22780.001

$ws.Range("H134").Value = 4994.5
$ws.Range("I134").Value = 4994.5
$ws.Range("K134").Value = 14983.5
$ws.Range("M134").Value = -12448.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 22349.691
$ws.Range("I14").Value = 22349.691
$ws.Range("K14").Value = 67049.073
$ws.Range("M14").Value = -66876.073

$ws.Range("H37").Value = 66000
$ws.Range("J37").Value = 66000
$ws.Range("L37").Value = 198000
$ws.Range("N37").Value = -198224

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 400000
$ws.Range("I5").Value = 400000
$ws.Range("K5").Value = 400000
$ws.Range("M5").Value = -399888

$ws.Range("H102").Value = 1552.2727
$ws.Range("I102").Value = 1657.8948
$ws.Range("J102").Value = 883.3333
$ws.Range("K102").Value = 1657.8948
$ws.Range("L102").Value = 883.3333
$ws.Range("M102").Value = -35.89480000000003
$ws.Range("N102").Value = -4127.3333

$ws.Range("H107").Value = 150.27272
$ws.Range("I107").Value = 95.25
$ws.Range("K107").Value = 95.25
$ws.Range("M107").Value = 1824.75

$ws.Range("H122").Value = 50186.668
$ws.Range("I122").Value = 51066
$ws.Range("K122").Value = 153198
$ws.Range("M122").Value = -150748

$ws.Range("H132").Value = 3233.5
$ws.Range("J132").Value = 2467
$ws.Range("L132").Value = 7401
$ws.Range("N132").Value = -12461

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 4998
$ws.Range("J12").Value = 4998
$ws.Range("L12").Value = 4998
$ws.Range("N12").Value = -5338

$ws.Range("H16").Value = 4514.778
$ws.Range("I16").Value = 4514.778
$ws.Range("K16").Value = 4514.778
$ws.Range("M16").Value = -4344.778

$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 5000
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6498

$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 25000
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -32488

$ws.Range("H95").Value = 39000
$ws.Range("J95").Value = 39000
$ws.Range("L95").Value = 39000
$ws.Range("N95").Value = -44492

$ws.Range("H122").Value = 4268.6665
$ws.Range("I122").Value = 4204.75
$ws.Range("J122").Value = 4396.5
$ws.Range("K122").Value = 12614.25
$ws.Range("L122").Value = 13189.5
$ws.Range("M122").Value = -10164.25
$ws.Range("N122").Value = -18089.5

$ws.Range("H136").Value = 4318.6665
$ws.Range("I136").Value = 2851.75
$ws.Range("K136").Value = 8555.25
$ws.Range("M136").Value = -6005.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 1429171.2
$ws.Range("I15").Value = 1429171.2
$ws.Range("K15").Value = 1429171.2
$ws.Range("M15").Value = -1428883.2

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H92").Value = 44999
$ws.Range("J92").Value = 44999
$ws.Range("L92").Value = 44999
$ws.Range("N92").Value = -49991

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H107").Value = 297.55554
$ws.Range("I107").Value = 298.2857
$ws.Range("K107").Value = 894.8571000000001
$ws.Range("M107").Value = 1025.1429

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 2419
$ws.Range("I122").Value = 1988.8334
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 5966.5002
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3516.5002
$ws.Range("N122").Value = -19900

$ws.Range("H126").Value = 4736.75
$ws.Range("I126").Value = 4649.1665
$ws.Range("K126").Value = 13947.4995
$ws.Range("M126").Value = -11477.4995

$ws.Range("H132").Value = 3826.4
$ws.Range("I132").Value = 3044
$ws.Range("K132").Value = 9132
$ws.Range("M132").Value = -6602

$ws.Range("H136").Value = 4991.25
$ws.Range("J136").Value = 4980
$ws.Range("L136").Value = 14940
$ws.Range("N136").Value = -20040
